# Updated railway-freight and inland-shipping data
# (Güterverkehr der Binnenschifffahrt, GENESIS-Tabelle 46321-0002)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab to match the GENESIS table id.
$ws.Name = "46321-0002"

# Row 359 (2020 / Februar): revised figures.
$ws.Cells.Item(359, 3).Value = 4178767
$ws.Cells.Item(359, 4).Value = 793
$ws.Cells.Item(359, 5).Value = 4247850
$ws.Cells.Item(359, 6).Value = 1144
$ws.Cells.Item(359, 7).Value = 6485123

# Row 360 (2020 / März): previously placeholder "..." cells now have
# published figures.
$ws.Cells.Item(360, 3).Value = 4782968
$ws.Cells.Item(360, 4).Value = 948
$ws.Cells.Item(360, 5).Value = 4163070
$ws.Cells.Item(360, 6).Value = 1118
$ws.Cells.Item(360, 7).Value = 7354357
$ws.Cells.Item(360, 8).Value = 1669
$ws.Cells.Item(360, 9).Value = 1059068
$ws.Cells.Item(360, 10).Value = 667

# Footer timestamp of the data export.
$ws.Cells.Item(373, 1).Value = "Stand: 06.07.2020 / 16:16:39"
